$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 30
$ws1.Range("F4").Value = 4795
$ws1.Range("F5").Value = 211
$ws1.Range("G5").Value = "不可售"
$ws1.Range("F6").Value = 155
$ws1.Range("F8").Value = 107
$ws1.Range("F9").Value = 93
$ws1.Range("F10").Value = 747
$ws1.Range("F11").Value = 222
$ws1.Range("F12").Value = 1167
$ws1.Range("F13").Value = 106
$ws1.Range("F15").Value = 182
$ws1.Range("F18").Value = 110
$ws1.Range("F20").Value = 6283
$ws1.Range("F21").Value = 39
$ws1.Range("F23").Value = 85
$ws1.Range("F26").Value = 3967
$ws1.Range("F27").Value = 402
$ws1.Range("F28").Value = 40
$ws1.Range("F29").Value = 13
$ws1.Range("F30").Value = 2568
$ws1.Range("F32").Value = 529
$ws1.Range("F34").Value = 284
$ws1.Range("F35").Value = 303
$ws1.Range("F36").Value = 369
$ws1.Range("F37").Value = 172
$ws1.Range("F39").Value = 1560
$ws1.Range("F40").Value = 960
$ws1.Range("F41").Value = 46
$ws1.Range("F42").Value = 64
$ws1.Range("F44").Value = 493
$ws1.Range("F47").Value = 583

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 30
$ws4.Range("F4").Value = 4795
$ws4.Range("F5").Value = 211
$ws4.Range("G5").Value = "不可售"
$ws4.Range("F6").Value = 155
$ws4.Range("F8").Value = 110
$ws4.Range("F9").Value = 107
$ws4.Range("F11").Value = 747
$ws4.Range("F12").Value = 222
$ws4.Range("F13").Value = 1167
$ws4.Range("F14").Value = 106
$ws4.Range("F15").Value = 259
$ws4.Range("F16").Value = 182
$ws4.Range("F17").Value = 82
$ws4.Range("F18").Value = 146
$ws4.Range("F19").Value = 110
$ws4.Range("F20").Value = 3952
$ws4.Range("F21").Value = 6283
$ws4.Range("F22").Value = 39
$ws4.Range("F23").Value = 0
$ws4.Range("F24").Value = 85
$ws4.Range("F26").Value = 48
$ws4.Range("F27").Value = 3967
$ws4.Range("F28").Value = 402
$ws4.Range("F29").Value = 40
$ws4.Range("F30").Value = 13
$ws4.Range("F31").Value = 2568
$ws4.Range("F33").Value = 529
$ws4.Range("F35").Value = 284
$ws4.Range("F36").Value = 303
$ws4.Range("F38").Value = 172
$ws4.Range("F40").Value = 1560
$ws4.Range("F41").Value = 960
$ws4.Range("F42").Value = 46
$ws4.Range("F43").Value = 64
$ws4.Range("F44").Value = 0
$ws4.Range("F46").Value = 479
$ws4.Range("F47").Value = 75
$ws4.Range("F48").Value = 583
